$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete the 50-problem list: the last three rows (N Queens, Minimum
# Window Substring, Egg Dropping) only had S.No / Title / Category filled
# in. Give them the same look as the rest of the finished table (copy the
# formatting down from the previous row) and fill in their Time
# Complexity / Optimised Time Complexity / Technique Used columns.

$ws.Range("A48:F48").Copy() | Out-Null
$ws.Range("A49:F51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D49").Value = "O(N!)"
$ws.Range("E49").Value = "O(N^2)"
$ws.Range("F49").Value = "backtracking"

$ws.Range("D50").Value = "O(n+m)"
$ws.Range("E50").Value = "O(n+m)"
$ws.Range("F50").Value = "sliding window + hashing"

$ws.Range("D51").Value = "O(k+n)"
$ws.Range("E51").Value = "O(k+n)"
$ws.Range("F51").Value = "binary search + 2d dynamic programming"

# Leave the selection where the author ended up after typing the last cell.
$ws.Range("F58").Select() | Out-Null
